$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2043256.8
$ws.Range("J17").Value = 2176441.8
$ws.Range("L17").Value = 6529325.399999999
$ws.Range("N17").Value = -6529661.399999999

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 95.75
$ws.Range("I33").Value = 95.75
$ws.Range("K33").Value = 95.75
$ws.Range("M33").Value = 133.25

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 182900.14
$ws.Range("J129").Value = 201156.95
$ws.Range("L129").Value = 603470.8500000001
$ws.Range("N129").Value = -613470.8500000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 4586.684
$ws.Range("I132").Value = 5436.467
$ws.Range("K132").Value = 16309.401
$ws.Range("M132").Value = -13779.401

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 10872483
$ws.Range("I135").Value = 409.11905
$ws.Range("J135").Value = 125029256
$ws.Range("K135").Value = 3682.07145
$ws.Range("L135").Value = 1125263304
$ws.Range("M135").Value = -1147.07145
$ws.Range("N135").Value = -1125268374

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1848.6061
$ws.Range("I137").Value = 1688.6154
$ws.Range("J137").Value = 2442.8572
$ws.Range("K137").Value = 5065.8462
$ws.Range("L137").Value = 7328.571599999999
$ws.Range("M137").Value = -2515.8462
$ws.Range("N137").Value = -12428.5716

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 13701510
$ws.Range("I138").Value = 41668116
$ws.Range("K138").Value = 125004348
$ws.Range("M138").Value = -124999208

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1454.0834
$ws.Range("I141").Value = 842.86206
$ws.Range("J141").Value = 3986.2856
$ws.Range("K141").Value = 2528.58618
$ws.Range("L141").Value = 11958.8568
$ws.Range("M141").Value = 2651.41382
$ws.Range("N141").Value = -22318.8568

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1943.73
$ws.Range("I32").Value = 1900.0215
$ws.Range("J32").Value = 2524.4285
$ws.Range("K32").Value = 1900.0215
$ws.Range("L32").Value = 2524.4285
$ws.Range("M32").Value = -1613.0215
$ws.Range("N32").Value = -3098.4285

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 564157.9
$ws.Range("I61").Value = 622009.0600000001
$ws.Range("J61").Value = 4929.3335
$ws.Range("K61").Value = 622009.0600000001
$ws.Range("L61").Value = 4929.3335
$ws.Range("M61").Value = -621797.0600000001
$ws.Range("N61").Value = -5353.3335

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 25001684
$ws.Range("I74").Value = 26317458
$ws.Range("J74").Value = 2000
$ws.Range("K74").Value = 26317458
$ws.Range("L74").Value = 2000
$ws.Range("M74").Value = -26316584
$ws.Range("N74").Value = -3748

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 25001684
$ws.Range("I77").Value = 26317458
$ws.Range("J77").Value = 2000
$ws.Range("K77").Value = 131587290
$ws.Range("L77").Value = 10000
$ws.Range("M77").Value = -131582922
$ws.Range("N77").Value = -18736

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 15278.108
$ws.Range("I132").Value = 1940.2307
$ws.Range("J132").Value = 46804
$ws.Range("K132").Value = 5820.6921
$ws.Range("L132").Value = 140412
$ws.Range("M132").Value = -3290.6921
$ws.Range("N132").Value = -145472

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 564157.9
$ws.Range("I136").Value = 622009.0600000001
$ws.Range("J136").Value = 4929.3335
$ws.Range("K136").Value = 1866027.18
$ws.Range("L136").Value = 14788.0005
$ws.Range("M136").Value = -1863477.18
$ws.Range("N136").Value = -19888.0005

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1516.5883
$ws.Range("I20").Value = 1800.3
$ws.Range("J20").Value = 1111.2858
$ws.Range("K20").Value = 1800.3
$ws.Range("L20").Value = 1111.2858
$ws.Range("M20").Value = -1553.3
$ws.Range("N20").Value = -1605.2858

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H61").Value = 17333.334
$ws.Range("J61").Value = 17333.334
$ws.Range("L61").Value = 17333.334
$ws.Range("N61").Value = -17959.334

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H110").Value = 20700.5
$ws.Range("J110").Value = 20700.5
$ws.Range("L110").Value = 20700.5
$ws.Range("N110").Value = -28880.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3639.611
$ws.Range("I134").Value = 3993.7856
$ws.Range("J134").Value = 2400
$ws.Range("K134").Value = 11981.3568
$ws.Range("L134").Value = 7200
$ws.Range("M134").Value = -9446.356800000001
$ws.Range("N134").Value = -12270

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4254.0513
$ws.Range("I31").Value = 2802.625
$ws.Range("K31").Value = 2802.625
$ws.Range("M31").Value = -2507.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H34").Value = 4254.0513
$ws.Range("I34").Value = 2802.625
$ws.Range("K34").Value = 2802.625
$ws.Range("M34").Value = -2600.625

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 8355.282999999999
$ws.Range("I58").Value = 743.08887
$ws.Range("J58").Value = 23925.682
$ws.Range("K58").Value = 743.08887
$ws.Range("L58").Value = 23925.682
$ws.Range("M58").Value = -540.08887
$ws.Range("N58").Value = -24331.682

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 21743000
$ws.Range("I99").Value = 3681.818
$ws.Range("K99").Value = 3681.818
$ws.Range("M99").Value = -2183.818

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 21743000
$ws.Range("I126").Value = 3681.818
$ws.Range("K126").Value = 11045.454
$ws.Range("M126").Value = -8575.454000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 1877.8085
$ws.Range("I132").Value = 1486.119
$ws.Range("J132").Value = 5168
$ws.Range("K132").Value = 4458.357
$ws.Range("L132").Value = 15504
$ws.Range("M132").Value = -1928.357
$ws.Range("N132").Value = -20564

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 1104.1666
$ws.Range("I134").Value = 1045.6842
$ws.Range("J134").Value = 1326.4
$ws.Range("K134").Value = 3137.0526
$ws.Range("L134").Value = 3979.2
$ws.Range("M134").Value = -602.0526
$ws.Range("N134").Value = -9049.200000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 8355.282999999999
$ws.Range("I136").Value = 743.08887
$ws.Range("J136").Value = 23925.682
$ws.Range("K136").Value = 2229.26661
$ws.Range("L136").Value = 71777.046
$ws.Range("M136").Value = 320.7333899999999
$ws.Range("N136").Value = -76877.046

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 3640.7693
$ws.Range("I3").Value = 1802.7273
$ws.Range("K3").Value = 5408.1819
$ws.Range("M3").Value = -5296.1819

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 767.59186
$ws.Range("I131").Value = 521.6667
$ws.Range("J131").Value = 783.63043
$ws.Range("K131").Value = 1565.0001
$ws.Range("L131").Value = 2350.89129
$ws.Range("M131").Value = 3474.9999
$ws.Range("N131").Value = -12430.89129

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H139").Value = 1309.5333
$ws.Range("I139").Value = 1086.6786
$ws.Range("J139").Value = 4429.5
$ws.Range("K139").Value = 3260.0358
$ws.Range("L139").Value = 13288.5
$ws.Range("M139").Value = 1879.9642
$ws.Range("N139").Value = -23568.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H140").Value = 2096
$ws.Range("I140").Value = 1405.2727
$ws.Range("K140").Value = 4215.8181
$ws.Range("M140").Value = 964.1818999999996

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 20835386
$ws.Range("I102").Value = 22728988
$ws.Range("K102").Value = 22728988
$ws.Range("M102").Value = -22727366

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H121").Value = 20158.5
$ws.Range("J121").Value = 20158.5
$ws.Range("L121").Value = 20158.5
$ws.Range("N121").Value = -23652.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4772.552
$ws.Range("I126").Value = 3705.8823
$ws.Range("J126").Value = 6283.6665
$ws.Range("K126").Value = 11117.6469
$ws.Range("L126").Value = 18850.9995
$ws.Range("M126").Value = -8647.6469
$ws.Range("N126").Value = -23790.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 30257.475
$ws.Range("I132").Value = 4679.8
$ws.Range("J132").Value = 126173.75
$ws.Range("K132").Value = 14039.4
$ws.Range("L132").Value = 378521.25
$ws.Range("M132").Value = -11509.4
$ws.Range("N132").Value = -383581.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 894208.9399999999
$ws.Range("I122").Value = 2453924.8
$ws.Range("J122").Value = 2942.7856
$ws.Range("K122").Value = 7361774.399999999
$ws.Range("L122").Value = 8828.356800000001
$ws.Range("M122").Value = -7359324.399999999
$ws.Range("N122").Value = -13728.3568

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 1254.8636
$ws.Range("I132").Value = 1313.7567
$ws.Range("J132").Value = 943.5714
$ws.Range("K132").Value = 3941.2701
$ws.Range("L132").Value = 2830.7142
$ws.Range("M132").Value = -1411.2701
$ws.Range("N132").Value = -7890.7142

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H136").Value = 1214.8572
$ws.Range("I136").Value = 1214.8572
$ws.Range("K136").Value = 3644.5716
$ws.Range("M136").Value = -1094.5716

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1463.6364
$ws.Range("I122").Value = 1488.8889
$ws.Range("J122").Value = 1350
$ws.Range("K122").Value = 4466.6667
$ws.Range("L122").Value = 4050
$ws.Range("M122").Value = -2016.6667
$ws.Range("N122").Value = -8950

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 1087.6923
$ws.Range("I126").Value = 1061.6666
$ws.Range("K126").Value = 3184.9998
$ws.Range("M126").Value = -714.9998000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 624.9016
$ws.Range("I132").Value = 515.14813
$ws.Range("J132").Value = 1471.5714
$ws.Range("K132").Value = 1545.44439
$ws.Range("L132").Value = 4414.7142
$ws.Range("M132").Value = 984.5556099999999
$ws.Range("N132").Value = -9474.7142

Write-Output "Applied all Typhon_Profits updates"
